$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add a zero-width "_GoBack" bookmark right before the run in the first
#    paragraph (the "Class Specification" heading), i.e. immediately after
#    <w:pPr> and before <w:r>.
#
#    A degenerate (zero-length) Range sitting exactly at document position 0
#    confuses bookmark placement, so we temporarily insert a throw-away
#    character at position 0, anchor the bookmark right after it (position
#    1, a perfectly ordinary mid-run position), and then delete the
#    throw-away character again. The bookmark correctly slides back down to
#    position 0 and serializes as tight bookmarkStart/bookmarkEnd pair.
# ---------------------------------------------------------------------------
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")
$anchor = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $anchor)
$d.Range(0, 1).Delete()

# ---------------------------------------------------------------------------
# 2. Table text clean-ups: drop w:proofErr spellStart/spellEnd wrappers and
#    merge runs that used to be split only because of those wrappers.
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)

# Row 1 / Col 2: "DatabaseInterface"
$cell = $tbl.Cell(1, 2)
$p = $cell.Range.Paragraphs(1)
$p.Range.InsertParagraphAfter()
$cellRange = $cell.Range
$newPara = $cellRange.Paragraphs(2)
$oldPara = $cellRange.Paragraphs(1)
$newPara.Style = $oldPara.Style
$newPara.Range.Text = "DatabaseInterface"
$oldPara.Range.Delete()

# Row 4 / Col 2: "CustomerController: 1" / "ReceivingController: 1"
$cell = $tbl.Cell(4, 2)
$cellRange = $cell.Range
$p1 = $cellRange.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$cellRange = $cell.Range
$newPara = $cellRange.Paragraphs(2)
$oldPara = $cellRange.Paragraphs(1)
$newPara.Style = $oldPara.Style
$newPara.Range.Text = "CustomerController: 1"
$oldPara.Range.Delete()

$cellRange = $cell.Range
$p2 = $cellRange.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$cellRange = $cell.Range
$newPara = $cellRange.Paragraphs(3)
$oldPara = $cellRange.Paragraphs(2)
$newPara.Style = $oldPara.Style
$newPara.Range.Text = "ReceivingController: 1"
$oldPara.Range.Delete()

# Row 5 / Col 2: "viewParts()" (leave the "Query database..." paragraph alone)
$cell = $tbl.Cell(5, 2)
$cellRange = $cell.Range
$p1 = $cellRange.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$cellRange = $cell.Range
$newPara = $cellRange.Paragraphs(2)
$oldPara = $cellRange.Paragraphs(1)
$newPara.Style = $oldPara.Style
$newPara.Range.Text = "viewParts()"
$oldPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Remove the stray "_GoBack" bookmark that used to sit at the end of the
#    "Revision and Date" cell (Row 10 / Col 2), now that it lives at the top
#    of the document instead.
# ---------------------------------------------------------------------------
$cell = $tbl.Cell(10, 2)
$cellRange = $cell.Range
$p1 = $cellRange.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$cellRange = $cell.Range
$newPara = $cellRange.Paragraphs(2)
$oldPara = $cellRange.Paragraphs(1)
$newPara.Style = $oldPara.Style
$newPara.Range.Text = "1.0; 4/7/2015"
$oldPara.Range.Delete()

Write-Output "done"
